$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 data values
$ws.Range("A2").Value = "Frank Warnakula"
$ws.Range("F2").Value = "sas"
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = "dd"

# Update the active selection to J2 (as reflected in the saved view state)
$ws.Range("J2").Select()
